$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Replace the small "time/size/education" table on Sheet2 with a full
# "Prime Minister / start_date / end_date" table.
$ws.Range("A1").Value = "Prime Minister"
$ws.Range("B1").Value = "start_date "
$ws.Range("C1").Value = "end_date"

$names = @(
  "Husyen Shaheed Suharwardy ",
  "Zulfikar Ali Bhutto",
  "Benazir Bhutto",
  "Benazir Bhutto",
  "Nawaz Sharif",
  "Nawaz Sharif",
  "Shahid Khaqan Abbasi",
  "Yousaf Raza Gillani",
  "Shehbaz Saharif",
  "Imran Khan"
)
$starts = @(22647, 28369, 31625, 35166, 36452, 43294, 43664, 43595, 44308, 45174)
$ends   = @(23012, 28949, 31717, 35394, 36445, 43796, 43918, 43614, 44832, 45265)

for ($i = 0; $i -lt $names.Count; $i++) {
  $r = $i + 2
  $ws.Cells.Item($r, 1).Value = $names[$i]
  $ws.Cells.Item($r, 2).Value = $starts[$i]
  $ws.Cells.Item($r, 2).NumberFormat = "m/d/yy"
  $ws.Cells.Item($r, 3).Value = $ends[$i]
  $ws.Cells.Item($r, 3).NumberFormat = "m/d/yy"
}

# Match column widths/best-fit for the date columns.
$ws.Columns("B:C").AutoFit()

# Portrait page setup for Sheet2.
$ws.PageSetup.Orientation = $excel.xlPortrait

# Make Sheet2 the active sheet/tab, with the given selection.
$ws.Activate()
$ws.Range("L15").Select()
